$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Efna4"
$ws.Range("C2").Value = "Epha4"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.828516
$ws.Range("H2").Value = 2.485548
$ws.Range("I2").Value = 0.4625620436231038
$ws.Range("J2").Value = 0.4821955800271095
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 7.727270999999999
$ws.Range("N2").Value = 23.181813
$ws.Range("O2").Value = 0.4492515923977304
$ws.Range("P2").Value = 0.4784711627054499
$ws.Range("Q2").Value = 6.402167659836
$ws.Range("R2").Value = 57.619508938524
$ws.Range("S2").Value = 0.2078067346804278
$ws.Range("T2").Value = 0.2307166798269999

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Efna4"
$ws.Range("C3").Value = "Epha4"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.828516
$ws.Range("H3").Value = 2.485548
$ws.Range("I3").Value = 0.4625620436231038
$ws.Range("J3").Value = 0.4821955800271095
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 6.292848333333333
$ws.Range("N3").Value = 18.878545
$ws.Range("O3").Value = 0.3658564756519351
$ws.Range("P3").Value = 0.3896519817642027
$ws.Range("Q3").Value = 5.21372552974
$ws.Range("R3").Value = 46.92352976766
$ws.Range("S3").Value = 0.1692313190503054
$ws.Range("T3").Value = 0.1878884633555024

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Efna4"
$ws.Range("C4").Value = "Epha4"
$ws.Range("D4").Value = "M1"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.828516
$ws.Range("H4").Value = 2.485548
$ws.Range("I4").Value = 0.4625620436231038
$ws.Range("J4").Value = 0.4821955800271095
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.021228
$ws.Range("N4").Value = 0.063684
$ws.Range("O4").Value = 0.001234163109255392
$ws.Range("P4").Value = 0.001314433755709007
$ws.Range("Q4").Value = 0.017587737648
$ws.Range("R4").Value = 0.158289638832
$ws.Range("S4").Value = 0.000570877009981418
$ws.Range("T4").Value = 0.0006338141472413165

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Efna4"
$ws.Range("C5").Value = "Epha4"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.828516
$ws.Range("H5").Value = 2.485548
$ws.Range("I5").Value = 0.4625620436231038
$ws.Range("J5").Value = 0.4821955800271095
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.007773333333333333
$ws.Range("N5").Value = 0.02332
$ws.Range("O5").Value = 0.0004519295852621654
$ws.Range("P5").Value = 0.0004813233336965963
$ws.Range("Q5").Value = 0.00644033104
$ws.Range("R5").Value = 0.05796297936
$ws.Range("S5").Value = 0.000209045472532609
$ws.Range("T5").Value = 0.0002320919840724122

# Row 6
$ws.Range("A6").Value = "ECs"
$ws.Range("B6").Value = "Efna4"
$ws.Range("C6").Value = "Epha4"
$ws.Range("D6").Value = "sCs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.828516
$ws.Range("H6").Value = 2.485548
$ws.Range("I6").Value = 0.4625620436231038
$ws.Range("J6").Value = 0.4821955800271095
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 3.151199
$ws.Range("N6").Value = 6.302398
$ws.Range("O6").Value = 0.183205839255817
$ws.Range("P6").Value = 0.1300810984409417
$ws.Range("Q6").Value = 2.610818790684
$ws.Range("R6").Value = 15.664912744104
$ws.Range("S6").Value = 0.08474406740985654
$ws.Range("T6").Value = 0.06272453071329342

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Efna4"
$ws.Range("C7").Value = "Epha4"
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.7438396666666667
$ws.Range("H7").Value = 2.231519
$ws.Range("I7").Value = 0.4152870872032183
$ws.Range("J7").Value = 0.4329140288365043
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 7.727270999999999
$ws.Range("N7").Value = 23.181813
$ws.Range("O7").Value = 0.4492515923977304
$ws.Range("P7").Value = 0.4784711627054499
$ws.Range("Q7").Value = 5.747850684883
$ws.Range("R7").Value = 51.730656163947
$ws.Range("S7").Value = 0.186568385228261
$ws.Range("T7").Value = 0.2071368787289029

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Efna4"
$ws.Range("C8").Value = "Epha4"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.7438396666666667
$ws.Range("H8").Value = 2.231519
$ws.Range("I8").Value = 0.4152870872032183
$ws.Range("J8").Value = 0.4329140288365043
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 6.292848333333333
$ws.Range("N8").Value = 18.878545
$ws.Range("O8").Value = 0.3658564756519351
$ws.Range("P8").Value = 0.3896519817642027
$ws.Range("Q8").Value = 4.680870206650555
$ws.Range("R8").Value = 42.127831859855
$ws.Range("S8").Value = 0.1519354701079273
$ws.Range("T8").Value = 0.1686858092696691

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Efna4"
$ws.Range("C9").Value = "Epha4"
$ws.Range("D9").Value = "M1"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.7438396666666667
$ws.Range("H9").Value = 2.231519
$ws.Range("I9").Value = 0.4152870872032183
$ws.Range("J9").Value = 0.4329140288365043
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.021228
$ws.Range("N9").Value = 0.063684
$ws.Range("O9").Value = 0.001234163109255392
$ws.Range("P9").Value = 0.001314433755709007
$ws.Range("Q9").Value = 0.015790228444
$ws.Range("R9").Value = 0.142112055996
$ws.Range("S9").Value = 0.0005125320027763391
$ws.Range("T9").Value = 0.0005690368128226835

# Row 10
$ws.Range("A10").Value = "FAPs"
$ws.Range("B10").Value = "Efna4"
$ws.Range("C10").Value = "Epha4"
$ws.Range("D10").Value = "M2"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.7438396666666667
$ws.Range("H10").Value = 2.231519
$ws.Range("I10").Value = 0.4152870872032183
$ws.Range("J10").Value = 0.4329140288365043
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.007773333333333333
$ws.Range("N10").Value = 0.02332
$ws.Range("O10").Value = 0.0004519295852621654
$ws.Range("P10").Value = 0.0004813233336965963
$ws.Range("Q10").Value = 0.005782113675555556
$ws.Range("R10").Value = 0.05203902308
$ws.Range("S10").Value = 0.0001876805210844832
$ws.Range("T10").Value = 0.0002083716235636106

# Row 11
$ws.Range("A11").Value = "FAPs"
$ws.Range("B11").Value = "Efna4"
$ws.Range("C11").Value = "Epha4"
$ws.Range("D11").Value = "sCs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.7438396666666667
$ws.Range("H11").Value = 2.231519
$ws.Range("I11").Value = 0.4152870872032183
$ws.Range("J11").Value = 0.4329140288365043
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 3.151199
$ws.Range("N11").Value = 6.302398
$ws.Range("O11").Value = 0.183205839255817
$ws.Range("P11").Value = 0.1300810984409417
$ws.Range("Q11").Value = 2.343986813760333
$ws.Range("R11").Value = 14.063920882562
$ws.Range("S11").Value = 0.07608301934316926
$ws.Range("T11").Value = 0.05631393240154599

# Row 12
$ws.Range("A12").Value = "sCs"
$ws.Range("B12").Value = "Efna4"
$ws.Range("C12").Value = "Epha4"
$ws.Range("D12").Value = "ECs"
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.21879
$ws.Range("H12").Value = 0.43758
$ws.Range("I12").Value = 0.1221508691736778
$ws.Range("J12").Value = 0.08489039113638626
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 7.727270999999999
$ws.Range("N12").Value = 23.181813
$ws.Range("O12").Value = 0.4492515923977304
$ws.Range("P12").Value = 0.4784711627054499
$ws.Range("Q12").Value = 1.69064962209
$ws.Range("R12").Value = 10.14389773254
$ws.Range("S12").Value = 0.05487647248904162
$ws.Range("T12").Value = 0.04061760414954715

# Row 13
$ws.Range("A13").Value = "sCs"
$ws.Range("B13").Value = "Efna4"
$ws.Range("C13").Value = "Epha4"
$ws.Range("D13").Value = "FAPs"
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.21879
$ws.Range("H13").Value = 0.43758
$ws.Range("I13").Value = 0.1221508691736778
$ws.Range("J13").Value = 0.08489039113638626
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 6.292848333333333
$ws.Range("N13").Value = 18.878545
$ws.Range("O13").Value = 0.3658564756519351
$ws.Range("P13").Value = 0.3896519817642027
$ws.Range("Q13").Value = 1.37681228685
$ws.Range("R13").Value = 8.2608737211
$ws.Range("S13").Value = 0.04468968649370238
$ws.Range("T13").Value = 0.03307770913903121

# Row 14
$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "Efna4"
$ws.Range("C14").Value = "Epha4"
$ws.Range("D14").Value = "M1"
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.21879
$ws.Range("H14").Value = 0.43758
$ws.Range("I14").Value = 0.1221508691736778
$ws.Range("J14").Value = 0.08489039113638626
$ws.Range("K14").Value = 1
$ws.Range("L14").Value = 0.3333333333333333
$ws.Range("M14").Value = 0.021228
$ws.Range("N14").Value = 0.063684
$ws.Range("O14").Value = 0.001234163109255392
$ws.Range("P14").Value = 0.001314433755709007
$ws.Range("Q14").Value = 0.00464447412
$ws.Range("R14").Value = 0.02786684472
$ws.Range("S14").Value = 0.0001507540964976349
$ws.Range("T14").Value = 0.0001115827956450068

# Row 15
$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "Efna4"
$ws.Range("C15").Value = "Epha4"
$ws.Range("D15").Value = "M2"
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.21879
$ws.Range("H15").Value = 0.43758
$ws.Range("I15").Value = 0.1221508691736778
$ws.Range("J15").Value = 0.08489039113638626
$ws.Range("K15").Value = 1
$ws.Range("L15").Value = 0.3333333333333333
$ws.Range("M15").Value = 0.007773333333333333
$ws.Range("N15").Value = 0.02332
$ws.Range("O15").Value = 0.0004519295852621654
$ws.Range("P15").Value = 0.0004813233336965963
$ws.Range("Q15").Value = 0.0017007276
$ws.Range("R15").Value = 0.0102043656
$ws.Range("S15").Value = 0.00005520359164507326
$ws.Range("T15").Value = 0.00004085972606057342

# Row 16
$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "Efna4"
$ws.Range("C16").Value = "Epha4"
$ws.Range("D16").Value = "sCs"
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.21879
$ws.Range("H16").Value = 0.43758
$ws.Range("I16").Value = 0.1221508691736778
$ws.Range("J16").Value = 0.08489039113638626
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 3.151199
$ws.Range("N16").Value = 6.302398
$ws.Range("O16").Value = 0.183205839255817
$ws.Range("P16").Value = 0.1300810984409417
$ws.Range("Q16").Value = 0.68945082921
$ws.Range("R16").Value = 2.75780331684
$ws.Range("S16").Value = 0.02237875250279115
$ws.Range("T16").Value = 0.01104263532610231
